$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-LogoInlineShape($range, [string]$newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $ishp = $range.InlineShapes.Item($i)
        # Round-trip through a floating Shape so we can write a new Name,
        # then convert back to an inline shape to keep <wp:inline> layout.
        $shp = $ishp.ConvertToShape()
        $shp.Name = $newName
        [void]$shp.ConvertToInlineShape()
    }
}

# Headers: BTec_Logo-Orange images renamed image1.jpg -> image2.jpg
Rename-LogoInlineShape $sec.Headers.Item(1).Range "image2.jpg"
Rename-LogoInlineShape $sec.Headers.Item(2).Range "image2.jpg"

# Footers: PearsonLogo images renamed image2.png -> image1.png
Rename-LogoInlineShape $sec.Footers.Item(1).Range "image1.png"
Rename-LogoInlineShape $sec.Footers.Item(2).Range "image1.png"
